# Update cryptos list with refreshed price/volume data (and a couple of
# re-ordered coin rows), matching the upstream "cryptos.xlsx" data refresh.
# Values that look like plain numbers (e.g. "1.00", "0.0606") are written
# with a leading apostrophe so Excel keeps them as text, exactly like the
# other price cells in this sheet (which are also stored as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.200.90'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '1.592.50'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''212.38'
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -0.36%  '
$ws.Range("D9").Value = '''0.0606'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").Value = '''19.02'
$ws.Range("E10").Value = '  -2.04%  '
$ws.Range("D11").Value = '''0.0847'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '1.817.14'
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").Value = '1.594.45'
$ws.Range("E13").Value = '  +0.24%  '
$ws.Range("D14").Value = '''4.01'
$ws.Range("E14").Value = '  -1.36%  '
$ws.Range("D15").Value = '''0.509'
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '26.199.20'
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '''7.37'
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''214.30'
$ws.Range("E20").Value = '  +1.69%  '
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").Value = '''4.25'
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("D23").Value = '''9.04'
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("D25").Value = '''144.89'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("D27").Value = '''6.95'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("E30").Value = '  -2.39%  '
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("E32").Value = '  -1.39%  '
$ws.Range("D33").Value = '1.426.14'
$ws.Range("E33").Value = '  +8.08%  '
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("E37").Value = '  -3.73%  '
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("D39").Value = '''0.824'
$ws.Range("E39").Value = '  +2.08%  '
$ws.Range("E40").Value = '  +4.53%  '
$ws.Range("D41").Value = '''0.999'
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = '''0.976'
$ws.Range("E42").Value = '  -10.08%  '
$ws.Range("D43").Value = '''0.767'
$ws.Range("E43").Value = '  +0.14%  '
$ws.Range("E44").Value = '  -0.18%  '
$ws.Range("D45").Value = '1.728.45'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").Value = '''61.04'
$ws.Range("E46").Value = '  -2.07%  '
$ws.Range("D47").Value = '''86.99'
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.49'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.0502'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '''0.0956'
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '''0.998'
$ws.Range("E51").Value = '  -0.23%  '
